$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 147.72728
$ws.Range("I11").Value = 147.72728
$ws.Range("K11").Value = 147.72728
$ws.Range("M11").Value = -7.727280000000007
$ws.Range("H40").Value = 2130.8076
$ws.Range("I40").Value = 1730.5385
$ws.Range("J40").Value = 2531.077
$ws.Range("K40").Value = 1730.5385
$ws.Range("L40").Value = 2531.077
$ws.Range("M40").Value = -1555.5385
$ws.Range("N40").Value = -2881.077
$ws.Range("H70").Value = 175892.25
$ws.Range("I70").Value = 2200
$ws.Range("J70").Value = 696969
$ws.Range("K70").Value = 6600
$ws.Range("L70").Value = 2090907
$ws.Range("M70").Value = -6330
$ws.Range("N70").Value = -2091447
$ws.Range("H73").Value = 175892.25
$ws.Range("I73").Value = 2200
$ws.Range("J73").Value = 696969
$ws.Range("K73").Value = 6600
$ws.Range("L73").Value = 2090907
$ws.Range("M73").Value = -5664
$ws.Range("N73").Value = -2092779
$ws.Range("H98").Value = 2035
$ws.Range("I98").Value = 2035
$ws.Range("K98").Value = 2035
$ws.Range("M98").Value = -537
$ws.Range("H100").Value = 750
$ws.Range("I100").Value = 750
$ws.Range("K100").Value = 750
$ws.Range("M100").Value = -209
$ws.Range("H122").Value = 2035
$ws.Range("I122").Value = 2035
$ws.Range("K122").Value = 6105
$ws.Range("M122").Value = -3655

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5366.8545
$ws.Range("I32").Value = 2742
$ws.Range("K32").Value = 2742
$ws.Range("M32").Value = -2455
$ws.Range("H45").Value = 3428.4285
$ws.Range("J45").Value = 3666.6667
$ws.Range("L45").Value = 3666.6667
$ws.Range("N45").Value = -4420.6667
$ws.Range("H61").Value = 1017.75
$ws.Range("I61").Value = 877.5714
$ws.Range("K61").Value = 877.5714
$ws.Range("M61").Value = -665.5714
$ws.Range("H132").Value = 3615.1538
$ws.Range("I132").Value = 2999.6667
$ws.Range("K132").Value = 8999.000100000001
$ws.Range("M132").Value = -6469.000100000001
$ws.Range("H136").Value = 1017.75
$ws.Range("I136").Value = 877.5714
$ws.Range("K136").Value = 2632.7142
$ws.Range("M136").Value = -82.71420000000035

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3916.6667
$ws.Range("I20").Value = 3916.6667
$ws.Range("K20").Value = 3916.6667
$ws.Range("M20").Value = -3669.6667
$ws.Range("H103").Value = 8090.5
$ws.Range("J103").Value = 8090.5
$ws.Range("L103").Value = 8090.5
$ws.Range("N103").Value = -10434.5
$ws.Range("H105").Value = 4647.8
$ws.Range("I105").Value = 7326.3335
$ws.Range("J105").Value = 3499.8572
$ws.Range("K105").Value = 7326.3335
$ws.Range("L105").Value = 3499.8572
$ws.Range("M105").Value = -5579.3335
$ws.Range("N105").Value = -6993.8572
$ws.Range("H134").Value = 1562.3334
$ws.Range("I134").Value = 1141.6
$ws.Range("K134").Value = 3424.8
$ws.Range("M134").Value = -889.7999999999997

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 1200
$ws.Range("I122").Value = 1200
$ws.Range("K122").Value = 3600
$ws.Range("M122").Value = -1150
$ws.Range("H125").Value = 44999
$ws.Range("J125").Value = 44999
$ws.Range("L125").Value = 44999
$ws.Range("N125").Value = -49919
$ws.Range("H132").Value = 1223.2
$ws.Range("I132").Value = 1127.25
$ws.Range("J132").Value = 1607
$ws.Range("K132").Value = 3381.75
$ws.Range("L132").Value = 4821
$ws.Range("M132").Value = -851.75
$ws.Range("N132").Value = -9881
$ws.Range("H134").Value = 3264.2727
$ws.Range("I134").Value = 2844.5715
$ws.Range("K134").Value = 8533.7145
$ws.Range("M134").Value = -5998.7145

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 510.94116
$ws.Range("I4").Value = 445.9655
$ws.Range("K4").Value = 1337.8965
$ws.Range("M4").Value = -1225.8965
$ws.Range("H9").Value = 349.5
$ws.Range("J9").Value = 349
$ws.Range("L9").Value = 1047
$ws.Range("N9").Value = -1495
$ws.Range("H36").Value = 450
$ws.Range("I36").Value = 450
$ws.Range("K36").Value = 1350
$ws.Range("M36").Value = -1181
$ws.Range("H120").Value = 14166.667
$ws.Range("I120").Value = 0
$ws.Range("J120").Value = 14166.667
$ws.Range("K120").Value = 0
$ws.Range("L120").Value = 42500.001
$ws.Range("M120").ClearContents()
$ws.Range("N120").Value = -52176.001
$ws.Range("H137").Value = 4327.9
$ws.Range("J137").Value = 3708.1667
$ws.Range("L137").Value = 11124.5001
$ws.Range("N137").Value = -21324.5001

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8500
$ws.Range("I70").Value = 8500
$ws.Range("K70").Value = 8500
$ws.Range("M70").Value = -8230
$ws.Range("H73").Value = 8500
$ws.Range("I73").Value = 8500
$ws.Range("K73").Value = 8500
$ws.Range("M73").Value = -7564
$ws.Range("H80").Value = 3528.04
$ws.Range("J80").Value = 4207.923
$ws.Range("L80").Value = 4207.923
$ws.Range("N80").Value = -6203.923
$ws.Range("H83").Value = 3528.04
$ws.Range("J83").Value = 4207.923
$ws.Range("L83").Value = 21039.615
$ws.Range("N83").Value = -31023.615
$ws.Range("H97").Value = 2327.923
$ws.Range("I97").Value = 539.6667
$ws.Range("J97").Value = 3860.7144
$ws.Range("K97").Value = 539.6667
$ws.Range("L97").Value = 3860.7144
$ws.Range("M97").Value = -43.66669999999999
$ws.Range("N97").Value = -4852.7144
$ws.Range("H102").Value = 2100.5789
$ws.Range("I102").Value = 1187.9231
$ws.Range("J102").Value = 4078
$ws.Range("K102").Value = 1187.9231
$ws.Range("L102").Value = 4078
$ws.Range("M102").Value = 434.0769
$ws.Range("N102").Value = -7322
$ws.Range("H126").Value = 4750.8
$ws.Range("I126").Value = 4713.75
$ws.Range("K126").Value = 14141.25
$ws.Range("M126").Value = -11671.25

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("I7").Value = 2000
$ws.Range("K7").Value = 2000
$ws.Range("M7").Value = -1888
$ws.Range("H25").Value = 19999.5
$ws.Range("I25").Value = 9999
$ws.Range("J25").Value = 30000
$ws.Range("K25").Value = 9999
$ws.Range("L25").Value = 30000
$ws.Range("M25").Value = -9769
$ws.Range("N25").Value = -30460
$ws.Range("H41").Value = 23665.666
$ws.Range("I41").Value = 22000
$ws.Range("J41").Value = 24498.5
$ws.Range("K41").Value = 22000
$ws.Range("L41").Value = 24498.5
$ws.Range("M41").Value = -21562
$ws.Range("N41").Value = -25374.5
$ws.Range("H100").Value = 1267.6666
$ws.Range("I100").Value = 1267.6666
$ws.Range("K100").Value = 1267.6666
$ws.Range("M100").Value = -726.6666
$ws.Range("H122").Value = 3165.5
$ws.Range("I122").Value = 2897.6667
$ws.Range("K122").Value = 8693.000100000001
$ws.Range("M122").Value = -6243.000100000001
$ws.Range("I126").Value = 2000
$ws.Range("K126").Value = 6000
$ws.Range("M126").Value = -3530

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1481.2632
$ws.Range("I132").Value = 1502.5883
$ws.Range("K132").Value = 4507.7649
$ws.Range("M132").Value = -1977.7649
